# Edit script for Pivots.xlsx
# - Updates the "Low" figure for column F (row 3) on the Povit sheet
# - Populates column K on the Povit sheet (a new trading-day column) by
#   copying the formatting/formulas from column J and filling in the
#   day's raw input values (Close/High/Low + a handful of manual pivot
#   reference values), letting Excel recompute every dependent formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Povit")

# --- Correct the existing "Low" value for 23-Jul (column F) ---
$ws.Range("F3").Value = 10138.6

# --- Copy column J's formatting (fill/border/number format) into column K
#     for every row that participates in the pivot calculations, so the
#     new column looks and behaves exactly like the others. ---
$ws.Range("J6:J46").Copy()
$ws.Range("K6:K46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Raw input values (Close / Low / High) for the new day in column K ---
$ws.Range("K2").Value = 10335.950000000001
$ws.Range("K3").Value = 10138.6
$ws.Range("K4").Value = 10234.65

# --- Pivot / support / resistance formulas for column K (mirrors column J) ---
$ws.Range("K6").Formula  = "=K8+K43"
$ws.Range("K7").Formula  = "=K11+K43"
$ws.Range("K8").Formula  = "=(2*K11)-K3"
$ws.Range("K10").Formula = "=K11+K13/2"
$ws.Range("K11").Formula = "=(K2+K3+K4)/3"
$ws.Range("K12").Formula = "=K11-K13/2"
$ws.Range("K13").Formula = "=ABS((K11-K46)*2)"
$ws.Range("K15").Formula = "=2*K11-K2"
$ws.Range("K16").Formula = "=K11-K43"
$ws.Range("K17").Formula = "=K15-K43"
$ws.Range("K19").Formula = "=(K2/K3)*K4"
$ws.Range("K20").Formula = "=K21+1.168*(K21-K22)"
$ws.Range("K21").Formula = "=K4+K44/2"
$ws.Range("K22").Formula = "=K4+K44/4"
$ws.Range("K23").Formula = "=K4+K44/6"
$ws.Range("K24").Formula = "=K4+K44/12"
$ws.Range("K25").Formula = "=K4"
$ws.Range("K26").Formula = "=K4-K44/12"
$ws.Range("K27").Formula = "=K4-K44/6"
$ws.Range("K28").Formula = "=K4-K44/4"
$ws.Range("K29").Formula = "=K4-K44/2"
$ws.Range("K30").Formula = "=K29-1.168*(K28-K29)"
$ws.Range("K31").Formula = "=K4-(K19-K4)"
$ws.Range("K37").Formula = "=K4"
$ws.Range("K43").Formula = "=ABS(K2-K3)"
$ws.Range("K44").Formula = "=K43*1.1"
$ws.Range("K45").Formula = "=(K2+K3)"
$ws.Range("K46").Formula = "=(K2+K3)/2"

# --- Manually recorded EW / Camarilla levels for the new day ---
$ws.Range("K35").Value = 10559
$ws.Range("K36").Value = 10335
$ws.Range("K38").Value = 10120
$ws.Range("K39").Value = 9951

# --- Recalculate everything and update the active selection on the sheet ---
$excel.CalculateFull()
$ws.Activate()
$ws.Range("K36").Select()
